$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- SamplesTab (row 3, column B): append "Order By" / LIMIT clause to existing query ---
$origSamples = $ws.Cells.Item(3, 2).Text
$ws.Range("B3").Value = $origSamples + "`nOrder By samp.sample_id_id ASC LIMIT 100"

# --- FilesTab (row 4, column B): reorder RETURN fields + append "Order By" / LIMIT clause ---
$newFiles = @'
MATCH (f:file)-->(parent)
WITH DISTINCT f, parent
MATCH (diag:diagnosis)-->(c)
OPTIONAL MATCH (f)-[*]->(samp:sample)
MATCH (f)-[*]->(c:case)<--(demo:demographic)
  WHERE diag.disease_term IN ['B Cell Lymphoma']

OPTIONAL MATCH (s:study)<--(c)<--(diag:diagnosis)<-[*]-(samp)
WITH
        f, parent, c, demo, diag, s, samp,
        ['Bytes', 'KB', 'MB', 'GB', 'TB'] AS units,
        toInteger(floor(log(f.file_size)/log(1024))) as i,
        2 as precision
WITH
        f, parent, c, demo, diag, s, samp,
        f.file_size /(1024^i) AS value, 
        10^precision AS factor,
        units[i] as unit
WITH    
        f, parent, c, demo, diag, s, samp, unit,
        round(factor * value)/factor AS size
RETURN 
        coalesce(f.file_name, '') AS `File Name`,
        coalesce(f.file_format, '') AS `Format`,
        coalesce(f.file_type, '') AS `File Type`,
       CASE size % 1 WHEN 0 THEN apoc.convert.toInteger(size)+' ' +unit ELSE size+' ' +unit END AS Size,
        coalesce(labels(parent)[0], '') AS `Association`,
        coalesce(f.file_description, '') AS `Description`,
        coalesce(samp.sample_id, '') AS `Sample ID`,
        coalesce(c.case_id, '') AS `Case ID`,
        coalesce(demo.breed,'') AS Breed ,
        coalesce(diag.disease_term,'') AS Diagnosis
Order By f.file_name ASC LIMIT 100
'@
$ws.Range("B4").Value = $newFiles

# --- StudyFilesTab (row 5, column B): append "Order By" / LIMIT clause to existing query ---
$origStudyFiles = $ws.Cells.Item(5, 2).Text
$ws.Range("B5").Value = $origStudyFiles + "`nOrder By f.file_name ASC LIMIT 100"

# --- CasesTab (row 2, column B): replace with rewritten query (reformatted + Order by/limit) ---
$newCases = @'
MATCH (s:study)<-[*]-(c:case)<--(demo:demographic)
MATCH (c)<--(diag:diagnosis)
  WHERE diag.disease_term IN ['B Cell Lymphoma']
OPTIONAL MATCH (samp:sample)-->(c)
OPTIONAL MATCH (co:cohort)<-[*]-(c)
WITH DISTINCT c, s, demo, diag, co, demo.patient_age_at_enrollment AS age, demo.weight as weight
RETURN  
coalesce(c.case_id, '') AS `Case ID` ,
        coalesce(s.clinical_study_designation, '') AS `Study Code` ,
        coalesce(s.clinical_study_type, '') AS  `Study Type`,
        coalesce(demo.breed, '') AS Breed ,
        coalesce(diag.disease_term, '') AS Diagnosis ,
        coalesce(diag.stage_of_disease, '') AS `Stage of Disease` ,
  coalesce(CASE age % 1 WHEN 0 THEN apoc.convert.toInteger(age) ELSE age END, '') AS Age,
       coalesce(demo.sex, '') AS Sex,
       coalesce(demo.neutered_indicator, '') AS `Neutered Status`,
coalesce(CASE weight % 1 WHEN 0 THEN apoc.convert.toInteger(weight) ELSE weight END, '') AS `Weight (kg)`,
       coalesce(diag.best_response, '') AS `Response to Treatment`,
       coalesce(co.cohort_description, '') AS `Cohort`
order by c.case_id asc
limit 100
'@
$ws.Range("B2").Value = $newCases

# --- Update selection to match the saved view (B2 instead of B8) ---
[void]$ws.Range("B2").Select()
